$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update totals in the summary block -----------------------------------
# VALOR MORA total
$ws.Range("E11").Value = 711200
# Cant. Periodos
$ws.Range("F13").Value = 14

# --- Fix "YAMILE REYES BAYONA" rows (16-18): periods sorted ascending ------
# and updated Salario Basico
$ws.Range("E16").Value = "2308"
$ws.Range("E17").Value = "2309"
$ws.Range("E18").Value = "2310"
$ws.Range("G16").Value = 1423500
$ws.Range("G17").Value = 1423500
$ws.Range("G18").Value = 1423500

# --- Re-sort "KEINER STIVEN AGUILAR JIMENEZ" rows (19-28) ascending -------
$ws.Range("E19").Value = "2410"
$ws.Range("E20").Value = "2411"
$ws.Range("E21").Value = "2412"
$ws.Range("E22").Value = "2501"
$ws.Range("E23").Value = "2502"
$ws.Range("E24").Value = "2503"
$ws.Range("E25").Value = "2504"
$ws.Range("E26").Value = "2505"
$ws.Range("E27").Value = "2506"
$ws.Range("E28").Value = "2507"

# --- Insert a new row for period 2508 (part 1 of new account statements) --
$ws.Rows.Item(29).Insert()

# Carry the previous "closing row" formatting (thicker bottom border) down
# to the freshly inserted row ...
$ws.Range("B28:J28").Copy()
$ws.Range("B29:J29").PasteSpecial(-4122)
# ... and restyle the old closing row (28) to the regular interior-row look
$ws.Range("B27:J27").Copy()
$ws.Range("B28:J28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "1065871889"
$ws.Range("D29").Value = "KEINER STIVEN AGUILAR JIMENEZ"
$ws.Range("E29").Value = "2508"
$ws.Range("F29").Value = 52000
$ws.Range("G29").Value = 1300000

Write-Host "done"
